$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 907
$ws.Range("F3").Value = 1477
$ws.Range("F4").Value = 1147
$ws.Range("F5").Value = 537
$ws.Range("F6").Value = 230
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = 696
$ws.Range("F9").Value = 286
$ws.Range("F13").Value = 167
$ws.Range("F14").Value = 3569
$ws.Range("F15").Value = 20
$ws.Range("F16").Value = 16
$ws.Range("F17").Value = 446
$ws.Range("F20").Value = 296
$ws.Range("F24").Value = 682
$ws.Range("F25").Value = 66
$ws.Range("F26").Value = 269
$ws.Range("F27").Value = 975
$ws.Range("F29").Value = 1630
$ws.Range("F30").Value = 358

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 246
$ws.Range("F6").Value = 31
$ws.Range("F7").Value = 244

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 115

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 907
$ws.Range("F4").Value = 1477
$ws.Range("F5").Value = 1147
$ws.Range("F8").Value = 115
$ws.Range("F9").Value = 537
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 696
$ws.Range("F14").Value = 286
$ws.Range("F18").Value = 167
$ws.Range("F19").Value = 3570
$ws.Range("F20").Value = 20
$ws.Range("F21").Value = 16
$ws.Range("F22").Value = 246
$ws.Range("F23").Value = 446
$ws.Range("F26").Value = 296
$ws.Range("F28").Value = 31
$ws.Range("F31").Value = 244
$ws.Range("F34").Value = 682
$ws.Range("F38").Value = 66
$ws.Range("F39").Value = 269
$ws.Range("F40").Value = 975
$ws.Range("F42").Value = 1630
$ws.Range("F43").Value = 358
